$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Тип помещения" (room/premises type) values in column C for rows 3-6
# to reflect the new example data (various premises types instead of
# repeating "Квартира").
$ws.Range("C3").Value = "Машиноместо"
$ws.Range("C4").Value = "Апартаменты"
$ws.Range("C5").Value = "Кладовая"
$ws.Range("C6").Value = "Коммерческое помещение"
